# Auto-generated edit script: applies numeric cell updates across 8 sheets
# as described by the commit diff (profit/price recompute for Rafflesia_Profits).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1298
$ws.Range("I12").Value = 1298
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1298
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1128
$ws.Range("N12").ClearContents()
$ws.Range("H40").Value = 1464.7693
$ws.Range("J40").Value = 1324.5
$ws.Range("L40").Value = 1324.5
$ws.Range("N40").Value = -1674.5
$ws.Range("H43").Value = 1728.7142
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 2020.2
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 2020.2
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -2158.2
$ws.Range("H53").Value = 619.8
$ws.Range("I53").Value = 327.7143
$ws.Range("J53").Value = 1301.3334
$ws.Range("K53").Value = 327.7143
$ws.Range("L53").Value = 1301.3334
$ws.Range("M53").Value = 309.2857
$ws.Range("N53").Value = -2575.3334
$ws.Range("H92").Value = 298.25
$ws.Range("I92").Value = 298.25
$ws.Range("K92").Value = 298.25
$ws.Range("M92").Value = 949.75
$ws.Range("H132").Value = 5225.75
$ws.Range("I132").Value = 5225.75
$ws.Range("K132").Value = 15677.25
$ws.Range("M132").Value = -13147.25
$ws.Range("H135").Value = 599.6667
$ws.Range("I135").Value = 399.5
$ws.Range("K135").Value = 3595.5
$ws.Range("M135").Value = -1060.5
$ws.Range("J137").Value = 2000
$ws.Range("L137").Value = 6000
$ws.Range("N137").Value = -11100
$ws.Range("H138").Value = 2276.7856
$ws.Range("J138").Value = 5280
$ws.Range("L138").Value = 15840
$ws.Range("N138").Value = -26120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12465
$ws.Range("I32").Value = 12465
$ws.Range("K32").Value = 12465
$ws.Range("M32").Value = -12178
$ws.Range("H74").Value = 4359.6
$ws.Range("I74").Value = 4199.5
$ws.Range("K74").Value = 4199.5
$ws.Range("M74").Value = -3325.5
$ws.Range("H77").Value = 4359.6
$ws.Range("I77").Value = 4199.5
$ws.Range("K77").Value = 20997.5
$ws.Range("M77").Value = -16629.5
$ws.Range("H96").Value = 33474
$ws.Range("J96").Value = 33474
$ws.Range("L96").Value = 33474
$ws.Range("N96").Value = -38966
$ws.Range("H97").Value = 903.1667
$ws.Range("I97").Value = 683.8
$ws.Range("K97").Value = 683.8
$ws.Range("M97").Value = -187.8
$ws.Range("H132").Value = 10221.4
$ws.Range("I132").Value = 8887.5
$ws.Range("K132").Value = 26662.5
$ws.Range("M132").Value = -24132.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 339
$ws.Range("I12").Value = 301.25
$ws.Range("J12").Value = 490
$ws.Range("K12").Value = 301.25
$ws.Range("L12").Value = 490
$ws.Range("M12").Value = -133.25
$ws.Range("N12").Value = -826
$ws.Range("H86").Value = 1875
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1875
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1875
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -4121
$ws.Range("H89").Value = 1875
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1875
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 9375
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -20607
$ws.Range("H94").Value = 1078.7778
$ws.Range("I94").Value = 903
$ws.Range("K94").Value = 903
$ws.Range("M94").Value = -452
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H134").Value = 5277.077
$ws.Range("I134").Value = 3849.1
$ws.Range("K134").Value = 11547.3
$ws.Range("M134").Value = -9012.299999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 138.75
$ws.Range("I7").Value = 185
$ws.Range("J7").Value = 61.666668
$ws.Range("K7").Value = 185
$ws.Range("L7").Value = 61.666668
$ws.Range("M7").Value = -72
$ws.Range("N7").Value = -287.666668
$ws.Range("H16").Value = 2250.3333
$ws.Range("I16").Value = 2547.25
$ws.Range("J16").Value = 1656.5
$ws.Range("K16").Value = 2547.25
$ws.Range("L16").Value = 1656.5
$ws.Range("M16").Value = -2260.25
$ws.Range("N16").Value = -2230.5
$ws.Range("H95").Value = 9126.556
$ws.Range("J95").Value = 9126.556
$ws.Range("L95").Value = 9126.556
$ws.Range("N95").Value = -14618.556
$ws.Range("H99").Value = 8000
$ws.Range("I99").Value = 8000
$ws.Range("K99").Value = 8000
$ws.Range("M99").Value = -6502
$ws.Range("H113").Value = 2250.3333
$ws.Range("I113").Value = 2547.25
$ws.Range("J113").Value = 1656.5
$ws.Range("K113").Value = 2547.25
$ws.Range("L113").Value = 1656.5
$ws.Range("M113").Value = -377.25
$ws.Range("N113").Value = -5996.5
$ws.Range("H126").Value = 8000
$ws.Range("I126").Value = 8000
$ws.Range("K126").Value = 24000
$ws.Range("M126").Value = -21530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 608.875
$ws.Range("I68").Value = 636.2
$ws.Range("J68").Value = 596.4545000000001
$ws.Range("K68").Value = 1908.6
$ws.Range("L68").Value = 1789.3635
$ws.Range("M68").Value = -1097.6
$ws.Range("N68").Value = -3411.3635
$ws.Range("H71").Value = 608.875
$ws.Range("I71").Value = 636.2
$ws.Range("J71").Value = 596.4545000000001
$ws.Range("K71").Value = 5725.8
$ws.Range("L71").Value = 5368.0905
$ws.Range("M71").Value = -1669.8
$ws.Range("N71").Value = -13480.0905
$ws.Range("H80").Value = 1798.8334
$ws.Range("I80").Value = 1000.5
$ws.Range("J80").Value = 2198
$ws.Range("K80").Value = 3001.5
$ws.Range("L80").Value = 6594
$ws.Range("M80").Value = -2065.5
$ws.Range("N80").Value = -8466
$ws.Range("H83").Value = 1798.8334
$ws.Range("I83").Value = 1000.5
$ws.Range("J83").Value = 2198
$ws.Range("K83").Value = 9004.5
$ws.Range("L83").Value = 19782
$ws.Range("M83").Value = -4324.5
$ws.Range("N83").Value = -29142
$ws.Range("H92").Value = 258.16666
$ws.Range("I92").Value = 166.66667
$ws.Range("J92").Value = 349.66666
$ws.Range("K92").Value = 500.00001
$ws.Range("L92").Value = 1048.99998
$ws.Range("M92").Value = 747.99999
$ws.Range("N92").Value = -3544.99998
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H132").Value = 497.5
$ws.Range("I132").Value = 497.5
$ws.Range("K132").Value = 4477.5
$ws.Range("M132").Value = -1947.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H123").Value = 58275
$ws.Range("J123").Value = 58275
$ws.Range("L123").Value = 58275
$ws.Range("N123").Value = -63175

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3300.25
$ws.Range("I7").Value = 2914.5715
$ws.Range("K7").Value = 2914.5715
$ws.Range("M7").Value = -2802.5715
$ws.Range("H22").Value = 2000
$ws.Range("H27").Value = 2000
$ws.Range("H53").Value = 3800
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H61").Value = 723
$ws.Range("I61").Value = 802.75
$ws.Range("J61").Value = 563.5
$ws.Range("K61").Value = 802.75
$ws.Range("L61").Value = 563.5
$ws.Range("M61").Value = -600.75
$ws.Range("N61").Value = -967.5
$ws.Range("H113").Value = 723
$ws.Range("I113").Value = 802.75
$ws.Range("J113").Value = 563.5
$ws.Range("K113").Value = 802.75
$ws.Range("L113").Value = 563.5
$ws.Range("M113").Value = 1367.25
$ws.Range("N113").Value = -4903.5
$ws.Range("H120").Value = 50698
$ws.Range("J120").Value = 50698
$ws.Range("L120").Value = 50698
$ws.Range("N120").Value = -60374
$ws.Range("H126").Value = 3300.25
$ws.Range("I126").Value = 2914.5715
$ws.Range("K126").Value = 8743.7145
$ws.Range("M126").Value = -6273.7145
$ws.Range("H132").Value = 1250
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -9560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7130.25
$ws.Range("I62").Value = 3299
$ws.Range("J62").Value = 8407.333000000001
$ws.Range("K62").Value = 3299
$ws.Range("L62").Value = 8407.333000000001
$ws.Range("M62").Value = -2675
$ws.Range("N62").Value = -9655.333000000001
$ws.Range("H65").Value = 7130.25
$ws.Range("I65").Value = 3299
$ws.Range("J65").Value = 8407.333000000001
$ws.Range("K65").Value = 16495
$ws.Range("L65").Value = 42036.665
$ws.Range("M65").Value = -13375
$ws.Range("N65").Value = -48276.665
$ws.Range("H107").Value = 250
$ws.Range("I107").Value = 250
$ws.Range("K107").Value = 750
$ws.Range("M107").Value = 1170
$ws.Range("H113").Value = 222
$ws.Range("I113").Value = 222
$ws.Range("K113").Value = 666
$ws.Range("M113").Value = 1504
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820
$ws.Range("H132").Value = 533.3333
$ws.Range("I132").Value = 533.3333
$ws.Range("K132").Value = 1599.9999
$ws.Range("M132").Value = 930.0001
$ws.Range("H136").Value = 933.3333
$ws.Range("I136").Value = 933.3333
$ws.Range("K136").Value = 2799.9999
$ws.Range("M136").Value = -249.9998999999998

Write-Host "Applied 251 cell updates and 7 cell clears across 8 sheets."
